$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so values like
# "1.001" or "241.77" are not re-interpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.949.72"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.894.26"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "0.8282"
$ws.Range("E5").Value = "  +8.04%  "
$ws.Range("D6").Value = "241.77"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.3228"
$ws.Range("E8").Value = "  +5.87%  "
$ws.Range("D9").Value = "26.56"
$ws.Range("E9").Value = "  +4.79%  "
$ws.Range("D10").Value = "0.07021"
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("D11").Value = "0.08034"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "0.7477"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "1.888.39"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "5.198"
$ws.Range("D15").Value = "92.35"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "29.955.56"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "14.03"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "5.895"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "245.09"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "0.000007762"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "2.146.98"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "6.906"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "0.1605"
$ws.Range("E25").Value = "  +25.23%  "
$ws.Range("D26").Value = "166.98"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "9.183"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "18.85"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D29").Value = "2.074"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").Value = "1.365"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").Value = "1.521"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").Value = "4.264"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "0.05643"
$ws.Range("E33").Value = "  +7.39%  "
$ws.Range("D34").Value = "4.075"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "1.274"
$ws.Range("E35").Value = "  +2.66%  "
$ws.Range("D36").Value = "0.7308"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D38").Value = "0.01910"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "2.781"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "0.4412"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "71.97"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "5.946"
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("D43").Value = "0.8424"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "7.585"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "100.71"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D50").Value = "2.042.99"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "36.01"
$ws.Range("E51").Value = "  -0.37%  "

# Row 48/49: EnergySwap and Maker swapped rank order with updated data
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "991.16"
$ws.Range("E48").Value = "  +9.24%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.672"
$ws.Range("E49").Value = "  -0.49%  "

# Clear the temporary text format so cells keep their original (default) style
$ws.Range("D2:E51").ClearFormats()
